$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 626451.2
$ws.Range("I2").Value = 853.3333
$ws.Range("J2").Value = 1001809.9
$ws.Range("K2").Value = 853.3333
$ws.Range("L2").Value = 1001809.9
$ws.Range("M2").Value = -740.3333
$ws.Range("N2").Value = -1002035.9

$ws.Range("H43").Value = 5185.2
$ws.Range("I43").Value = 4981.75
$ws.Range("K43").Value = 4981.75
$ws.Range("M43").Value = -4912.75

$ws.Range("H64").Value = 3695.25
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 3651.7144
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 3651.7144
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -4147.7144

$ws.Range("H67").Value = 3695.25
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 3651.7144
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 3651.7144
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -5367.7144

$ws.Range("H137").Value = 2187.6924
$ws.Range("I137").Value = 2352.4443
$ws.Range("J137").Value = 1817
$ws.Range("K137").Value = 7057.3329
$ws.Range("L137").Value = 5451
$ws.Range("M137").Value = -4507.3329
$ws.Range("N137").Value = -10551


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3728
$ws.Range("I2").Value = 3775.2222
$ws.Range("J2").Value = 3643
$ws.Range("K2").Value = 3775.2222
$ws.Range("L2").Value = 3643
$ws.Range("M2").Value = -3662.2222
$ws.Range("N2").Value = -3869

$ws.Range("H45").Value = 3662.6667
$ws.Range("I45").Value = 2998.6667
$ws.Range("K45").Value = 2998.6667
$ws.Range("M45").Value = -2621.6667

$ws.Range("H61").Value = 15002897
$ws.Range("I61").Value = 19094052
$ws.Range("K61").Value = 19094052
$ws.Range("M61").Value = -19093840

$ws.Range("H63").Value = 4309.778
$ws.Range("J63").Value = 4697.3335
$ws.Range("L63").Value = 4697.3335
$ws.Range("N63").Value = -6069.3335

$ws.Range("H66").Value = 4309.778
$ws.Range("J66").Value = 4697.3335
$ws.Range("L66").Value = 23486.6675
$ws.Range("N66").Value = -30350.6675

$ws.Range("H74").Value = 2557.95
$ws.Range("I74").Value = 1509
$ws.Range("J74").Value = 4131.375
$ws.Range("K74").Value = 1509
$ws.Range("L74").Value = 4131.375
$ws.Range("M74").Value = -635
$ws.Range("N74").Value = -5879.375

$ws.Range("H77").Value = 2557.95
$ws.Range("I77").Value = 1509
$ws.Range("J77").Value = 4131.375
$ws.Range("K77").Value = 7545
$ws.Range("L77").Value = 20656.875
$ws.Range("M77").Value = -3177
$ws.Range("N77").Value = -29392.875

$ws.Range("H101").Value = 100000
$ws.Range("J101").Value = 100000
$ws.Range("L101").Value = 100000
$ws.Range("N101").Value = -106490

$ws.Range("H110").Value = 1391.3684
$ws.Range("I110").Value = 495.93332
$ws.Range("K110").Value = 495.93332
$ws.Range("M110").Value = 1549.06668

$ws.Range("H116").Value = 3728
$ws.Range("I116").Value = 3775.2222
$ws.Range("J116").Value = 3643
$ws.Range("K116").Value = 3775.2222
$ws.Range("L116").Value = 3643
$ws.Range("M116").Value = -1481.2222
$ws.Range("N116").Value = -8231

$ws.Range("H122").Value = 3826.5625
$ws.Range("I122").Value = 3941.6667
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 11825.0001
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -9375.000100000001
$ws.Range("N122").Value = -11200

$ws.Range("H132").Value = 4713.625
$ws.Range("I132").Value = 3968.7222
$ws.Range("J132").Value = 6948.3335
$ws.Range("K132").Value = 11906.1666
$ws.Range("L132").Value = 20845.0005
$ws.Range("M132").Value = -9376.1666
$ws.Range("N132").Value = -25905.0005

$ws.Range("H136").Value = 15002897
$ws.Range("I136").Value = 19094052
$ws.Range("K136").Value = 57282156
$ws.Range("M136").Value = -57279606


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3728
$ws.Range("I3").Value = 3775.2222
$ws.Range("J3").Value = 3643
$ws.Range("K3").Value = 3775.2222
$ws.Range("L3").Value = 3643
$ws.Range("M3").Value = -3661.2222
$ws.Range("N3").Value = -3871

$ws.Range("H86").Value = 60338.332
$ws.Range("I86").Value = 97245.71000000001
$ws.Range("J86").Value = 8668
$ws.Range("K86").Value = 97245.71000000001
$ws.Range("L86").Value = 8668
$ws.Range("M86").Value = -96122.71000000001
$ws.Range("N86").Value = -10914

$ws.Range("H89").Value = 60338.332
$ws.Range("I89").Value = 97245.71000000001
$ws.Range("J89").Value = 8668
$ws.Range("K89").Value = 486228.55
$ws.Range("L89").Value = 43340
$ws.Range("M89").Value = -480612.55
$ws.Range("N89").Value = -54572


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28574608
$ws.Range("I31").Value = 52634744
$ws.Range("K31").Value = 52634744
$ws.Range("M31").Value = -52634449

$ws.Range("H34").Value = 28574608
$ws.Range("I34").Value = 52634744
$ws.Range("K34").Value = 52634744
$ws.Range("M34").Value = -52634542

$ws.Range("H58").Value = 2493.8262
$ws.Range("I58").Value = 2331.6667
$ws.Range("J58").Value = 2797.875
$ws.Range("K58").Value = 2331.6667
$ws.Range("L58").Value = 2797.875
$ws.Range("M58").Value = -2128.6667
$ws.Range("N58").Value = -3203.875

$ws.Range("H134").Value = 3234.7727
$ws.Range("I134").Value = 3060.6316
$ws.Range("K134").Value = 9181.8948
$ws.Range("M134").Value = -6646.8948

$ws.Range("H136").Value = 2493.8262
$ws.Range("I136").Value = 2331.6667
$ws.Range("J136").Value = 2797.875
$ws.Range("K136").Value = 6995.000100000001
$ws.Range("L136").Value = 8393.625
$ws.Range("M136").Value = -4445.000100000001


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 72199.39999999999
$ws.Range("I99").Value = 29999
$ws.Range("J99").Value = 82749.5
$ws.Range("K99").Value = 29999
$ws.Range("L99").Value = 82749.5
$ws.Range("M99").Value = -27753
$ws.Range("N99").Value = -87241.5


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 71433530
$ws.Range("I61").Value = 125000936
$ws.Range("K61").Value = 125000936
$ws.Range("M61").Value = -125000734

$ws.Range("H104").Value = 35198.332
$ws.Range("J104").Value = 35198.332
$ws.Range("L104").Value = 35198.332
$ws.Range("N104").Value = -42186.332

$ws.Range("H113").Value = 71433530
$ws.Range("I113").Value = 125000936
$ws.Range("K113").Value = 125000936
$ws.Range("M113").Value = -124998766

$ws.Range("H136").Value = 3200.476
$ws.Range("I136").Value = 2853.2104
$ws.Range("K136").Value = 8559.6312
$ws.Range("M136").Value = -6009.6312


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1281.5454
$ws.Range("I81").Value = 1281.5454
$ws.Range("K81").Value = 2563.0908
$ws.Range("M81").Value = -1502.0908

$ws.Range("H84").Value = 1281.5454
$ws.Range("I84").Value = 1281.5454
$ws.Range("K84").Value = 12815.454
$ws.Range("M84").Value = -7511.454

$ws.Range("H107").Value = 5449
$ws.Range("I107").Value = 4726.2
$ws.Range("J107").Value = 5639.2104
$ws.Range("K107").Value = 14178.6
$ws.Range("L107").Value = 16917.6312
$ws.Range("M107").Value = -12258.6
$ws.Range("N107").Value = -20757.6312

$ws.Range("H136").Value = 14147.926
$ws.Range("I136").Value = 14119.8
$ws.Range("J136").Value = 14499.5
$ws.Range("K136").Value = 42359.39999999999
$ws.Range("L136").Value = 43498.5
$ws.Range("M136").Value = -39809.39999999999
$ws.Range("N136").Value = -48598.5

